$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log row (row 4) mirroring the existing knowledge-log entries
$ws.Range("A4").Value = "2025-08-26 18:02:57"
$ws.Range("B4").Value = "Climate Change"
$ws.Range("C4").Value = "General"
$ws.Range("D4").Value = 530
$ws.Range("E4").Value = 1570
$ws.Range("F4").Value = "SUCCESS"
$ws.Range("G4").Value = "Generated successfully"
